# Update scripts with new TPM values for the Nlgn3-Nrxn1 LR-pair sheet.
# Rows 2-7 are refreshed with recalculated TPM-based figures (and the "Target
# cluster" column is reassigned), and three new rows (8-10) are appended so the
# sheet now contains the full 3x3 Sending-cluster x Target-cluster matrix
# (ECs, FAPs, MuSCs) for the Nlgn3 -> Nrxn1 ligand-receptor pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ A="ECs"; B="Nlgn3"; C="Nrxn1"; D="ECs"; E=1; F=0.3333333333333333; G=0.007630333333333334; H=0.022891; I=0.005544871895800688; J=0.005544871895800688; K=2; L=0.6666666666666666; M=0.006517333333333333; N=0.019552; O=0.004697037253424763; P=0.004697037253424762; Q=0.00004972942577777778; R=0.000447564832; S=0.00002604446986004382; T=0.00002604446986004382 },
  @{ A="ECs"; B="Nlgn3"; C="Nrxn1"; D="FAPs"; E=1; F=0.3333333333333333; G=0.007630333333333334; H=0.022891; I=0.005544871895800688; J=0.005544871895800688; K=1; L=0.3333333333333333; M=0.01189366666666667; N=0.035681; O=0.008571756661182945; P=0.008571756661182945; Q=0.00009075264122222221; R=0.000816773771; S=0.00004752929260823565; T=0.00004752929260823565 },
  @{ A="ECs"; B="Nlgn3"; C="Nrxn1"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.007630333333333334; H=0.022891; I=0.005544871895800688; J=0.005544871895800688; K=3; L=1; M=1.369130333333333; N=4.107391; O=0.9867312060853923; P=0.9867312060853922; Q=0.01044692082011111; R=0.094022287381; S=0.005471298133332409; T=0.005471298133332408 },
  @{ A="FAPs"; B="Nlgn3"; C="Nrxn1"; D="ECs"; E=3; F=1; G=0.4678513333333334; H=1.403554; I=0.3399819636031033; J=0.3399819636031033; K=2; L=0.6666666666666666; M=0.006517333333333333; N=0.019552; O=0.004697037253424763; P=0.004697037253424762; Q=0.003049143089777778; R=0.027442287808; S=0.001596907948536278; T=0.001596907948536278 },
  @{ A="FAPs"; B="Nlgn3"; C="Nrxn1"; D="FAPs"; E=3; F=1; G=0.4678513333333334; H=1.403554; I=0.3399819636031033; J=0.3399819636031033; K=1; L=0.3333333333333333; M=0.01189366666666667; N=0.035681; O=0.008571756661182945; P=0.008571756661182945; Q=0.005564467808222223; R=0.050080210274; S=0.002914242661196959; T=0.002914242661196959 },
  @{ A="FAPs"; B="Nlgn3"; C="Nrxn1"; D="MuSCs"; E=3; F=1; G=0.4678513333333334; H=1.403554; I=0.3399819636031033; J=0.3399819636031033; K=3; L=1; M=1.369130333333333; N=4.107391; O=0.9867312060853923; P=0.9867312060853922; Q=0.6405494519571112; R=5.764945067614001; S=0.3354708129933701; T=0.3354708129933701 },
  @{ A="MuSCs"; B="Nlgn3"; C="Nrxn1"; D="ECs"; E=3; F=1; G=0.9006246666666667; H=2.701874; I=0.6544731645010959; J=0.6544731645010959; K=2; L=0.6666666666666666; M=0.006517333333333333; N=0.019552; O=0.004697037253424763; P=0.004697037253424762; Q=0.005869671160888889; R=0.052827040448; S=0.003074084835028441; T=0.00307408483502844 },
  @{ A="MuSCs"; B="Nlgn3"; C="Nrxn1"; D="FAPs"; E=3; F=1; G=0.9006246666666667; H=2.701874; I=0.6544731645010959; J=0.6544731645010959; K=1; L=0.3333333333333333; M=0.01189366666666667; N=0.035681; O=0.008571756661182945; P=0.008571756661182945; Q=0.01071172957711111; R=0.096405566194; S=0.00560998470737775; T=0.00560998470737775 },
  @{ A="MuSCs"; B="Nlgn3"; C="Nrxn1"; D="MuSCs"; E=3; F=1; G=0.9006246666666667; H=2.701874; I=0.6544731645010959; J=0.6544731645010959; K=3; L=1; M=1.369130333333333; N=4.107391; O=0.9867312060853923; P=0.9867312060853922; Q=1.233072550081556; R=11.097652950734; S=0.6457890949586897; T=0.6457890949586897 }
)

$rowNum = 2
foreach ($row in $rows) {
  foreach ($col in $row.Keys) {
    $ws.Range("$col$rowNum").Value = $row[$col]
  }
  $rowNum++
}
